$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "7+85="
$t.Cell(1,2).Range.Text = "19+62="
$t.Cell(1,3).Range.Text = "46-17="
$t.Cell(1,4).Range.Text = "13-8="
$t.Cell(1,5).Range.Text = "33-8="
$t.Cell(2,1).Range.Text = "54+18="
$t.Cell(2,2).Range.Text = "66+7="
$t.Cell(2,3).Range.Text = "8+77="
$t.Cell(2,4).Range.Text = "87-49="
$t.Cell(2,5).Range.Text = "9+49="
$t.Cell(3,1).Range.Text = "57-8="
$t.Cell(3,2).Range.Text = "71-52="
$t.Cell(3,3).Range.Text = "64+29="
$t.Cell(3,4).Range.Text = "17+44="
$t.Cell(3,5).Range.Text = "9+39="
$t.Cell(4,1).Range.Text = "45-39="
$t.Cell(4,2).Range.Text = "49+26="
$t.Cell(4,3).Range.Text = "21-13="
$t.Cell(4,4).Range.Text = "82-38="
$t.Cell(4,5).Range.Text = "39+55="
$t.Cell(5,1).Range.Text = "93-18="
$t.Cell(5,2).Range.Text = "3+78="
$t.Cell(5,3).Range.Text = "74-16="
$t.Cell(5,4).Range.Text = "41-15="
$t.Cell(5,5).Range.Text = "64-16="
$t.Cell(6,1).Range.Text = "81-57="
$t.Cell(6,2).Range.Text = "8+18="
$t.Cell(6,3).Range.Text = "44+8="
$t.Cell(6,4).Range.Text = "52-38="
$t.Cell(6,5).Range.Text = "68+7="
$t.Cell(7,1).Range.Text = "47+35="
$t.Cell(7,2).Range.Text = "66+17="
$t.Cell(7,3).Range.Text = "82-8="
$t.Cell(7,4).Range.Text = "26+9="
$t.Cell(7,5).Range.Text = "57+4="
$t.Cell(8,1).Range.Text = "46+15="
$t.Cell(8,2).Range.Text = "56+9="
$t.Cell(8,3).Range.Text = "18+55="
$t.Cell(8,4).Range.Text = "95-67="
$t.Cell(8,5).Range.Text = "8+7="
$t.Cell(9,1).Range.Text = "25+69="
$t.Cell(9,2).Range.Text = "25+67="
$t.Cell(9,3).Range.Text = "66+17="
$t.Cell(9,4).Range.Text = "81-68="
$t.Cell(9,5).Range.Text = "71-7="
$t.Cell(10,1).Range.Text = "18+43="
$t.Cell(10,2).Range.Text = "59+32="
$t.Cell(10,3).Range.Text = "9+34="
$t.Cell(10,4).Range.Text = "73-34="
$t.Cell(10,5).Range.Text = "39+48="
$t.Cell(11,1).Range.Text = "72-19="
$t.Cell(11,2).Range.Text = "19+23="
$t.Cell(11,3).Range.Text = "68-49="
$t.Cell(11,4).Range.Text = "26+39="
$t.Cell(11,5).Range.Text = "57+38="
$t.Cell(12,1).Range.Text = "47-18="
$t.Cell(12,2).Range.Text = "78-49="
$t.Cell(12,3).Range.Text = "35+37="
$t.Cell(12,4).Range.Text = "64+8="
$t.Cell(12,5).Range.Text = "34-8="
$t.Cell(13,1).Range.Text = "62-4="
$t.Cell(13,2).Range.Text = "18+64="
$t.Cell(13,3).Range.Text = "95-49="
$t.Cell(13,4).Range.Text = "79+4="
$t.Cell(13,5).Range.Text = "49+43="
$t.Cell(14,1).Range.Text = "69+29="
$t.Cell(14,2).Range.Text = "84-49="
$t.Cell(14,3).Range.Text = "37+37="
$t.Cell(14,4).Range.Text = "7+86="
$t.Cell(14,5).Range.Text = "69+14="
$t.Cell(15,1).Range.Text = "60-23="
$t.Cell(15,2).Range.Text = "3+89="
$t.Cell(15,3).Range.Text = "50-18="
$t.Cell(15,4).Range.Text = "66-47="
$t.Cell(15,5).Range.Text = "44-28="
$t.Cell(16,1).Range.Text = "42-38="
$t.Cell(16,2).Range.Text = "6+66="
$t.Cell(16,3).Range.Text = "22+49="
$t.Cell(16,4).Range.Text = "68+19="
$t.Cell(16,5).Range.Text = "30-14="
$t.Cell(17,1).Range.Text = "83-65="
$t.Cell(17,2).Range.Text = "12-8="
$t.Cell(17,3).Range.Text = "4+7="
$t.Cell(17,4).Range.Text = "5+28="
$t.Cell(17,5).Range.Text = "93-6="
$t.Cell(18,1).Range.Text = "57+28="
$t.Cell(18,2).Range.Text = "59+37="
$t.Cell(18,3).Range.Text = "90-31="
$t.Cell(18,4).Range.Text = "58+29="
$t.Cell(18,5).Range.Text = "19+5="
$t.Cell(19,1).Range.Text = "47+15="
$t.Cell(19,2).Range.Text = "25+67="
$t.Cell(19,3).Range.Text = "92-18="
$t.Cell(19,4).Range.Text = "52-14="
$t.Cell(19,5).Range.Text = "71-53="
$t.Cell(20,1).Range.Text = "7+9="
$t.Cell(20,2).Range.Text = "19+28="
$t.Cell(20,3).Range.Text = "29+38="
$t.Cell(20,4).Range.Text = "23-7="
$t.Cell(20,5).Range.Text = "38+49="
